$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40, shifting existing rows 40..125 down to 41..126.
$ws.Rows.Item(40).EntireRow.Insert()

# Populate the newly inserted row 40 with the new data record.
$ws.Range("A40").Value = 10
$ws.Range("B40").Value = "Vega Modelo de Temuco"
$ws.Range("C40").Value = "La Araucanía"
$ws.Range("D40").Value = 44720
$ws.Range("E40").Value = 9
$ws.Range("F40").Value = "Fruta"
$ws.Range("G40").Value = 100104
$ws.Range("H40").Value = "Frutos de pepita"
$ws.Range("I40").Value = 100104001
$ws.Range("J40").Value = "Granada"
$ws.Range("K40").Value = "Sin especificar"
$ws.Range("L40").Value = "Primera"
$ws.Range("M40").Value = 75
$ws.Range("N40").Value = 15000
$ws.Range("O40").Value = 16000
$ws.Range("P40").Value = 15467
$ws.Range("Q40").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R40").Value = "Provincia de Limarí"
$ws.Range("S40").Value = 1031
$ws.Range("T40").Value = 15
